# Update "想去人数" (want-to-go count) values in the F column of several sheets
# to the values captured at the latest data scrape (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 517
$ws.Range("F6").Value = 932
$ws.Range("F7").Value = 169
$ws.Range("F8").Value = 965
$ws.Range("F9").Value = 753
$ws.Range("F10").Value = 205
$ws.Range("F13").Value = 790
$ws.Range("F17").Value = 1307
$ws.Range("F21").Value = 2811
$ws.Range("F22").Value = 1325
$ws.Range("F23").Value = 662
$ws.Range("F27").Value = 977
$ws.Range("F28").Value = 324
$ws.Range("F29").Value = 1715
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 1346

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 513

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 517
$ws.Range("F8").Value = 513
$ws.Range("F9").Value = 513
$ws.Range("F13").Value = 932
$ws.Range("F14").Value = 169
$ws.Range("F16").Value = 965
$ws.Range("F17").Value = 753
$ws.Range("F18").Value = 205
$ws.Range("F26").Value = 790
$ws.Range("F30").Value = 1307
$ws.Range("F34").Value = 2811
$ws.Range("F35").Value = 1325
$ws.Range("F36").Value = 662
$ws.Range("F42").Value = 977
$ws.Range("F43").Value = 324
$ws.Range("F44").Value = 1715
$ws.Range("F46").Value = 4
$ws.Range("F47").Value = 1346
